$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 used to be "Task 6" -- replace it in place with the new "Uni Break" task.
# Copy the date format from an existing date cell first (so the new cell
# reuses the workbook's existing date style instead of creating a new one),
# then write the values.
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4122)
$ws.Cells.Item(7, 1).Value = "Uni Break"
$ws.Cells.Item(7, 2).Value = (Get-Date -Year 2023 -Month 11 -Day 19).Date
$ws.Cells.Item(7, 3).Value = 90

# Row 8 used to be "Task 7" -- replace it in place with the new "Prototype Research" task.
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(8, 2).PasteSpecial(-4122)
$ws.Cells.Item(8, 1).Value = "Prototype Research"
$ws.Cells.Item(8, 2).Value = (Get-Date -Year 2024 -Month 2 -Day 21).Date
$ws.Cells.Item(8, 3).Value = 20

# "Article Research" (row 5) days-to-complete bumped from 64 to 200
$ws.Cells.Item(5, 3).Value = 200

# Update the saved selection to match the new editing position
$ws.Range("C16").Select()
